# Auto-generated edit script: updates cryptos list price/volume columns
# (commit: "Updated cryptos list on Mon Sep 16 14:47:53 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.801.59'
$ws.Range('E2').Value = '  -3.89%  '
$ws.Range('D3').Value = '2.280.55'
$ws.Range('E3').Value = '  -5.30%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'542.20"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.87%  '
$ws.Range('D6').Value = "'130.74"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.34%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = "'0.568"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.20%  '
$ws.Range('D9').Value = '2.279.90'
$ws.Range('E9').Value = '  -5.18%  '
$ws.Range('D10').Value = "'0.0999"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -5.25%  '
$ws.Range('E11').Value = '  -3.13%  '
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('E13').Value = '  -5.47%  '
$ws.Range('D14').Value = "'23.39"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -5.33%  '
$ws.Range('D15').Value = '2.684.19'
$ws.Range('E15').Value = '  -5.44%  '
$ws.Range('D16').Value = '57.819.05'
$ws.Range('E16').Value = '  -3.74%  '
$ws.Range('E17').Value = '  -4.65%  '
$ws.Range('D18').Value = '2.216.34'
$ws.Range('E18').Value = '  -6.42%  '
$ws.Range('D19').Value = "'10.53"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -6.12%  '
$ws.Range('D20').Value = "'4.24"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -6.07%  '
$ws.Range('D21').Value = "'312.03"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.27%  '
$ws.Range('D22').Value = "'6.39"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -6.03%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = "'62.56"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.41%  '
$ws.Range('E25').Value = '  -3.36%  '
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').Value = "'7.94"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -6.99%  '
$ws.Range('E28').Value = '  -7.65%  '
$ws.Range('E29').Value = '  -4.17%  '
$ws.Range('D30').Value = "'169.86"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.47%  '
$ws.Range('D31').Value = '0.0₃0715'
$ws.Range('E31').Value = '  -6.82%  '
$ws.Range('D32').Value = "'1.09"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.94%  '
$ws.Range('D33').Value = "'5.71"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -6.55%  '
$ws.Range('D34').Value = "'0.378"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -5.89%  '
$ws.Range('D35').Value = "'0.998"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').Value = "'17.63"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -4.12%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').Value = "'1.23"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -8.09%  '
$ws.Range('D39').Value = "'3.89"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -6.90%  '
$ws.Range('D40').Value = "'37.85"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.77%  '
$ws.Range('E41').Value = '  -7.33%  '
$ws.Range('D42').Value = "'286.81"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -11.61%  '
$ws.Range('D43').Value = "'139.20"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -6.35%  '
$ws.Range('E44').Value = '  -4.97%  '
$ws.Range('D45').Value = "'0.0947"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.21%  '
$ws.Range('D46').Value = "'0.0497"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.57%  '
$ws.Range('D47').Value = "'0.550"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.43%  '
$ws.Range('D48').Value = "'18.16"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -8.87%  '
$ws.Range('D49').Value = "'0.0211"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -4.80%  '
$ws.Range('D50').Value = "'10.95"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('D51').Value = "'16.41"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.25%  '
